# Update the Traceability Matrix: rename "Fulfilled by" column to "Tested by"
# and replace its contents with the actual testing method for each requirement.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header + values for column D ("Tested by")
$testedBy = @{
    1  = "Tested by"
    2  = "Pressing floor buttons simultaneously"
    3  = "Pressing any floor button"
    4  = "Pressing floor button, then pressing elevator button"
    5  = "N/A - part of GUI"
    6  = "N/A - part of GUI"
    7  = "Pressing open/close doors button"
    8  = "Pressing open/close doors button"
    9  = "Pressing the help button"
    10 = "Pressing any floor button, integrated into movement"
    11 = "Testing with floor button then elevator button"
    12 = "N/A - part of GUI"
    13 = "Pressing the help button"
    14 = "Pressing the door obstacle button"
    15 = "Pressing the fire button"
    16 = "Setting the weight of the passenger >300 in the elevator.h class"
    17 = "Pressing the power outage button"
}

foreach ($row in 1..17) {
    $ws.Cells.Item($row, 4).Value = $testedBy[$row]
}

$ws.Range("B31").Select()
